$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "38.150.72"
$ws.Range("E2").Value = "  +2.94%  "
$ws.Range("D3").Value = "2.059.74"
$ws.Range("E3").Value = "  +2.43%  "
$ws.Range("E4").Value = "  -0.19%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "230.27"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +1.70%  "
$ws.Range("E6").Value = "  +1.52%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "58.19"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  +6.21%  "
$ws.Range("E9").Value = "  +2.20%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0806"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  +2.47%  "
$ws.Range("E11").Value = "  -1.20%  "
$ws.Range("D12").Value = "2.363.56"
$ws.Range("E12").Value = "  +2.46%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "14.61"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  +3.17%  "
$ws.Range("E14").Value = "  +2.06%  "
$ws.Range("E15").Value = "  +1.98%  "
$ws.Range("E16").Value = "  +3.40%  "
$ws.Range("D17").Value = "2.063.47"
$ws.Range("E17").Value = "  +2.96%  "
$ws.Range("D18").Value = "38.073.74"
$ws.Range("E18").Value = "  +2.96%  "
$ws.Range("E19").Value = "  +1.21%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "69.76"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  +1.38%  "
$ws.Range("E21").Value = "  +1.61%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "224.89"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  +0.67%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.999"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  -0.15%  "
$ws.Range("E24").Value = "  +1.20%  "
$ws.Range("E25").Value = "  +2.86%  "
$ws.Range("E26").Value = "  +2.05%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "166.03"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  +0.06%  "
$ws.Range("E28").Value = "  +8.04%  "
$ws.Range("E30").Value = "  +1.76%  "
$ws.Range("E31").Value = "  +1.54%  "
$ws.Range("E32").Value = "  +1.14%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.61"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  +4.46%  "
$ws.Range("E34").Value = "  +0.97%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.98"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  +7.00%  "
$ws.Range("E36").Value = "  +1.84%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "6.09"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  +14.14%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.33"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  +5.63%  "
$ws.Range("E39").Value = "  +0.01%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "98.42"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  +3.83%  "
$ws.Range("E41").Value = "  +1.01%  "
$ws.Range("D42").Value = "1.484.87"
$ws.Range("E42").Value = "  +0.47%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "16.89"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  +2.04%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0943"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  +2.30%  "
$ws.Range("E45").Value = "  +3.72%  "
$ws.Range("E46").Value = "  -0.16%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "4.10"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  +16.89%  "
$ws.Range("E48").Value = "  +1.27%  "
$ws.Range("E49").Value = "  +1.93%  "
$ws.Range("E50").Value = "  -1.03%  "
$ws.Range("D51").Value = "2.252.93"
$ws.Range("E51").Value = "  +2.69%  "
